$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-4 and add new rows 5-7 with recalculated NATMI values
# (new "ECs" cell-type cluster added to the analysis)

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Efna5"
$ws.Range("C2").Value = "Epha4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.666083666666667
$ws.Range("H2").Value = 4.998251
$ws.Range("I2").Value = 0.6125276070882968
$ws.Range("J2").Value = 0.6125276070882968
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 2.856403666666667
$ws.Range("N2").Value = 8.569211
$ws.Range("O2").Value = 0.235832554697756
$ws.Range("P2").Value = 0.235832554697756
$ws.Range("Q2").Value = 4.759007494440111
$ws.Range("R2").Value = 42.831067449961
$ws.Range("S2").Value = 0.1444539504025364
$ws.Range("T2").Value = 0.1444539504025363

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Efna5"
$ws.Range("C3").Value = "Epha4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.666083666666667
$ws.Range("H3").Value = 4.998251
$ws.Range("I3").Value = 0.6125276070882968
$ws.Range("J3").Value = 0.6125276070882968
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 7.000300666666668
$ws.Range("N3").Value = 21.000902
$ws.Range("O3").Value = 0.5779641054021444
$ws.Range("P3").Value = 0.5779641054021444
$ws.Range("Q3").Value = 11.66308660248911
$ws.Range("R3").Value = 104.967779422402
$ws.Range("S3").Value = 0.3540189704649037
$ws.Range("T3").Value = 0.3540189704649037

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Efna5"
$ws.Range("C4").Value = "Epha4"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.666083666666667
$ws.Range("H4").Value = 4.998251
$ws.Range("I4").Value = 0.6125276070882968
$ws.Range("J4").Value = 0.6125276070882968
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.255294666666666
$ws.Range("N4").Value = 6.765884
$ws.Range("O4").Value = 0.1862033399000996
$ws.Range("P4").Value = 0.1862033399000996
$ws.Range("Q4").Value = 3.757509607653777
$ws.Range("R4").Value = 33.817586468884
$ws.Range("S4").Value = 0.1140546862208568
$ws.Range("T4").Value = 0.1140546862208568

# Row 5
$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Efna5"
$ws.Range("C5").Value = "Epha4"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.053930333333333
$ws.Range("H5").Value = 3.161791
$ws.Range("I5").Value = 0.3874723929117032
$ws.Range("J5").Value = 0.3874723929117031
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 2.856403666666667
$ws.Range("N5").Value = 8.569211
$ws.Range("O5").Value = 0.235832554697756
$ws.Range("P5").Value = 0.235832554697756
$ws.Range("Q5").Value = 3.010450468544556
$ws.Range("R5").Value = 27.094054216901
$ws.Range("S5").Value = 0.09137860429521966
$ws.Range("T5").Value = 0.09137860429521963

# Row 6
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Efna5"
$ws.Range("C6").Value = "Epha4"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.053930333333333
$ws.Range("H6").Value = 3.161791
$ws.Range("I6").Value = 0.3874723929117032
$ws.Range("J6").Value = 0.3874723929117031
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 7.000300666666668
$ws.Range("N6").Value = 21.000902
$ws.Range("O6").Value = 0.5779641054021444
$ws.Range("P6").Value = 0.5779641054021444
$ws.Range("Q6").Value = 7.377829215053557
$ws.Range("R6").Value = 66.40046293548201
$ws.Range("S6").Value = 0.2239451349372407
$ws.Range("T6").Value = 0.2239451349372407

# Row 7
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Efna5"
$ws.Range("C7").Value = "Epha4"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.053930333333333
$ws.Range("H7").Value = 3.161791
$ws.Range("I7").Value = 0.3874723929117032
$ws.Range("J7").Value = 0.3874723929117031
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 2.255294666666666
$ws.Range("N7").Value = 6.765884
$ws.Range("O7").Value = 0.1862033399000996
$ws.Range("P7").Value = 0.1862033399000996
$ws.Range("Q7").Value = 2.376923459804889
$ws.Range("R7").Value = 21.392311138244
$ws.Range("S7").Value = 0.07214865367924281
$ws.Range("T7").Value = 0.07214865367924278
